$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.585.27"
$ws.Range("E2").Value = "  +3.59%  "

$ws.Range("D3").Value = "1.694.91"
$ws.Range("E3").Value = "  +1.96%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.12"
$ws.Range("E5").Value = "  +1.80%  "

$ws.Range("E6").Value = "  +0.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3939"
$ws.Range("E7").Value = "  +1.25%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4018"
$ws.Range("E8").Value = "  +1.85%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.521"
$ws.Range("E9").Value = "  +6.34%  "

$ws.Range("E10").Value = "  +0.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.64"
$ws.Range("E11").Value = "  +8.02%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08769"
$ws.Range("E12").Value = "  +1.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.233"
$ws.Range("E13").Value = "  +8.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.23"
$ws.Range("E14").Value = "  +2.58%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001324"
$ws.Range("E15").Value = "  +0.45%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.588"
$ws.Range("E16").Value = "  +5.08%  "

$ws.Range("D17").Value = "1.695.08"
$ws.Range("E17").Value = "  +1.70%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "100.08"
$ws.Range("E18").Value = "  +0.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07070"
$ws.Range("E19").Value = "  +4.27%  "

$ws.Range("E20").Value = "  +3.27%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.873"
$ws.Range("E21").Value = "  +3.75%  "

$ws.Range("E22").Value = "  +0.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.05"
$ws.Range("E23").Value = "  +1.58%  "

$ws.Range("D24").Value = "24.586.09"
$ws.Range("E24").Value = "  +3.62%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.026"
$ws.Range("E25").Value = "  +8.37%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.314"
$ws.Range("E26").Value = "  -0.46%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.42"
$ws.Range("E27").Value = "  +3.22%  "

$ws.Range("E28").Value = "  +0.71%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.248"
$ws.Range("E29").Value = "  +1.65%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.80"
$ws.Range("E30").Value = "  +4.10%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.433"
$ws.Range("E31").Value = "  +14.61%  "

$ws.Range("D32").Value = "1.884.52"
$ws.Range("E32").Value = "  +1.74%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.108"
$ws.Range("E33").Value = "  -1.14%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08531"
$ws.Range("E34").Value = "  -0.20%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.248"
$ws.Range("E35").Value = "  +10.16%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.40"
$ws.Range("E36").Value = "  +9.82%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.955"
$ws.Range("E37").Value = "  +0.50%  "

$ws.Range("E38").Value = "  +2.94%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.51"
$ws.Range("E39").Value = "  +0.75%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02745"
$ws.Range("E40").Value = "  +9.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.09062"
$ws.Range("E41").Value = "  +3.20%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.464"
$ws.Range("E42").Value = "  +1.29%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7727"
$ws.Range("E43").Value = "  +2.56%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7221"
$ws.Range("E44").Value = "  +3.01%  "

$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.547"
$ws.Range("E45").Value = "  +5.89%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "15.46"
$ws.Range("E46").Value = "  +4.48%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.203"
$ws.Range("E47").Value = "  +2.43%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.353"
$ws.Range("E48").Value = "  +13.55%  "

$ws.Range("E49").Value = "  +0.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "141.27"
$ws.Range("E50").Value = "  +2.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.08027"
$ws.Range("E51").Value = "  +3.46%  "
